$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sales rows recorded by the entry form. Columns F (Price) and G (Date)
# hold values that *look* numeric/date-like but must be stored as literal
# text (matching the source app's export), so those are entered with a
# leading apostrophe to force text entry, then the cell style is
# re-pointed at a plain text cell's style so the number/date auto-format
# doesn't stick to it.
$rows = @(
    @{ Row=3; A="c161e675-a9a0-49f6-957b-3fc3eedc5b99"; B="ooredoo"; C="SIM Card"; D="DIMA";  E=2;  F="1500";  G="2024-09-10"; H="20:38:49"; Reset=$false },
    @{ Row=4; A="b51f35fd-d5b3-45e0-807d-cb95078d64b5"; B="Master";  C="Other";    D="Other";  E=10; F="1000";  G="2024-09-10"; H="20:39:55"; Reset=$false },
    @{ Row=5; A="8b6e27cf-3bd7-4a7b-a085-5f9d1805f7b2"; B="Master";  C="Other";    D="Other";  E=10; F="1000";  G="2024-09-10"; H="21:34:40"; Reset=$false },
    @{ Row=6; A="67bec3f9-8d9b-4a02-9c4e-06b70a131d29"; B="LDNIO";   C="Cable";    D="Type C"; E=1;  F="400.0"; G="2024-09-10"; H="21:54:37"; Reset=$false },
    @{ Row=7; A="c40f6909-09f5-4ded-85c2-c5c153f7502b"; B="OOREDOO"; C="SIM Card"; D="DIMA";   E=1;  F="1000";  G="2024-09-10"; H="22:01:02"; Reset=$true  }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = "'" + $r.F
    $ws.Cells.Item($row, 7).Value = "'" + $r.G
    $ws.Cells.Item($row, 8).Value = $r.H

    # Re-base the quote-prefixed text cells onto a normal text cell's style
    # so they keep rendering/storing as plain text without tagging along
    # an extra quote-prefix style variant.
    $ws.Cells.Item($row, 6).Style = $ws.Cells.Item($row, 1).Style
    $ws.Cells.Item($row, 7).Style = $ws.Cells.Item($row, 1).Style

    if ($r.Reset) {
        # The form resets formatting back to the workbook default when the
        # user switches fields, so the last row picks up the plain/default
        # style instead of inheriting the sheet's usual cell style.
        $ws.Cells.Item($row, 1).Style = "Normal"
        $ws.Cells.Item($row, 2).Style = "Normal"
        $ws.Cells.Item($row, 3).Style = "Normal"
        $ws.Cells.Item($row, 4).Style = "Normal"
        $ws.Cells.Item($row, 5).Style = "Normal"
        $ws.Cells.Item($row, 6).Style = "Normal"
        $ws.Cells.Item($row, 7).Style = "Normal"
        $ws.Cells.Item($row, 8).Style = "Normal"
    }
}
